$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header for column B
$ws.Range("B1").Value = "Change to total individual income tax burden (`$ million)"

# Update the district aggregate values in column B
$ws.Range("B2").Value = 466
$ws.Range("B3").Value = 636
$ws.Range("B4").Value = 669
$ws.Range("B5").Value = 543
$ws.Range("B6").Value = 326
$ws.Range("B7").Value = 708
$ws.Range("B8").Value = 524
$ws.Range("B9").Value = 363
$ws.Range("B10").Value = 607
$ws.Range("B11").Value = 253
$ws.Range("B12").Value = 2105

# Remove column C entirely (percentage column no longer used)
$ws.Range("C1:C12").EntireColumn.Delete()
